# Insert a new row at row 15 (pushing GBDT..XGBoost rows down by one),
# and populate it with the new DeepCNN entry.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row before current row 15 (GBDT / subsample ...)
$ws.Rows.Item(15).Insert()

# Populate the new row 15 with DeepCNN in column A and empty strings B:F
$ws.Cells.Item(15, 1).Value = "DeepCNN"
$ws.Cells.Item(15, 2).Value = ""
$ws.Cells.Item(15, 3).Value = ""
$ws.Cells.Item(15, 4).Value = ""
$ws.Cells.Item(15, 5).Value = ""
$ws.Cells.Item(15, 6).Value = ""
